# Apply edit: replace existing row 2 data set and add new rows 2-4 and 6
# (old row 2 content is preserved but moved to row 5, with an updated date_of_lab)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "2024-04-17_00:00:00.000_IMTX_Conference_Note_91586"
$ws.Range("B2").Value = "'56.21"
$ws.Range("C2").Value = "<3.08"
$ws.Range("D2").Value = ">115.57"
$ws.Range("E2").Value = "'2024-02-22"
$ws.Range("F2").Value = "['Kappa Free Light Chain 0.76 - 6.83 mg/dL 56.21 (H)', 'Lambda Free Light Chain 0.68 - 4.58 mg/dL <3.08 (L)', 'Kappa/Lambda FLC Ratio 0.66 - 2.22  >115.57 (H)']"
$ws.Range("G2").Value = "{
  `"kappa_flc`": `"56.21`",
  `"lambda_flc`": `"<3.08`",
  `"kappa_lambda_ratio`": `">115.57`",
  `"date_of_lab`": `"2024-02-22`",
  `"evidence_sentences`": [
    `"Kappa Free Light Chain 0.76 - 6.83 mg/dL 56.21 (H)`",
    `"Lambda Free Light Chain 0.68 - 4.58 mg/dL <3.08 (L)`",
    `"Kappa/Lambda FLC Ratio 0.66 - 2.22  >115.57 (H)`"
  ],
  `"source_document`": `"2024-04-17_00:00:00.000_IMTX_Conference_Note_91586`"
}"

# Row 3
$ws.Range("A3").Value = "2024-04-17_00:00:00.000_IMTX_Conference_Note_91586"
$ws.Range("B3").Value = "'203.94"
$ws.Range("C3").Value = "<0.15"
$ws.Range("D3").Value = ">1456.71"
$ws.Range("E3").Value = "'2024-02-08"
$ws.Range("F3").Value = "['Kappa Free Light Chain 203.94 (H) 0.76 - 6.83 mg/dL', 'Lambda Free Light Chain <0.15 (L) 0.68 - 4.58 mg/dL', 'Kappa/Lambda FLC Ratio >1456.71 (H) 0.66 - 2.22']"
$ws.Range("G3").Value = "{
  `"kappa_flc`": `"203.94`",
  `"lambda_flc`": `"<0.15`",
  `"kappa_lambda_ratio`": `">1456.71`",
  `"date_of_lab`": `"2024-02-08`",
  `"evidence_sentences`": [
    `"Kappa Free Light Chain 203.94 (H) 0.76 - 6.83 mg/dL`",
    `"Lambda Free Light Chain <0.15 (L) 0.68 - 4.58 mg/dL`",
    `"Kappa/Lambda FLC Ratio >1456.71 (H) 0.66 - 2.22`"
  ],
  `"source_document`": `"2024-04-17_00:00:00.000_IMTX_Conference_Note_91586`"
}"

# Row 4
$ws.Range("A4").Value = "2024-02-08_00:00:00.000_Progress_Notes_91427"
$ws.Range("B4").Value = "'0.08"
$ws.Range("C4").Value = "<0.15"
$ws.Range("D4").Value = ">0.57"
$ws.Range("E4").Value = "'2024-06-13"
$ws.Range("F4").Value = "['Kappa Free Light Chain 0.08 (L) 0.76 - 6.83 mg/dL', 'Lambda Free Light Chain <0.15 (L) 0.68 - 4.58 mg/dL', 'Kappa/Lambda FLC Ratio >0.57 (L) 0.66 - 2.22']"
$ws.Range("G4").Value = "{
  `"kappa_flc`": `"0.08`",
  `"lambda_flc`": `"<0.15`",
  `"kappa_lambda_ratio`": `">0.57`",
  `"date_of_lab`": `"2024-06-13`",
  `"evidence_sentences`": [
    `"Kappa Free Light Chain 0.08 (L) 0.76 - 6.83 mg/dL`",
    `"Lambda Free Light Chain <0.15 (L) 0.68 - 4.58 mg/dL`",
    `"Kappa/Lambda FLC Ratio >0.57 (L) 0.66 - 2.22`"
  ],
  `"source_document`": `"2024-02-08_00:00:00.000_Progress_Notes_91427`"
}"

# Row 5
$ws.Range("A5").Value = "2024-06-25_00:00:00.000_Progress_Notes_91596"
$ws.Range("B5").Value = "<0.06 mg/dL"
$ws.Range("C5").Value = "<1.61 mg/dL"
$ws.Range("E5").Value = "'2024-06-25"
$ws.Range("F5").Value = "['Labs from 4/8/2024: Kappa <0.06 mg/dL, Lambda <1.61 mg/dL, SPEP with M-spike 0.3 g/dL, IgG kappa']"
$ws.Range("G5").Value = "{
  `"kappa_flc`": `"<0.06 mg/dL`",
  `"lambda_flc`": `"<1.61 mg/dL`",
  `"kappa_lambda_ratio`": null,
  `"date_of_lab`": `"2024-06-25`",
  `"evidence_sentences`": [
    `"Labs from 4/8/2024: Kappa <0.06 mg/dL, Lambda <1.61 mg/dL, SPEP with M-spike 0.3 g/dL, IgG kappa`"
  ],
  `"source_document`": `"2024-06-25_00:00:00.000_Progress_Notes_91596`"
}"

# Row 6
$ws.Range("A6").Value = "2024-04-17_00:00:00.000_Progress_Notes_91591"
$ws.Range("B6").Value = "'16.18"
$ws.Range("C6").Value = "<0.15"
$ws.Range("D6").Value = ">115.57"
$ws.Range("E6").Value = "'2024-04-17"
$ws.Range("F6").Value = "['2/22/24: KFLC 16.18, LFLC <0.15, kappa/lambda ratio >115.57. SPEP with IgG kappa, M-spike  1.0; previous monoclonal kappa free light chain not detected.']"
$ws.Range("G6").Value = "{
  `"kappa_flc`": `"16.18`",
  `"lambda_flc`": `"<0.15`",
  `"kappa_lambda_ratio`": `">115.57`",
  `"date_of_lab`": `"2024-04-17`",
  `"evidence_sentences`": [
    `"2/22/24: KFLC 16.18, LFLC <0.15, kappa/lambda ratio >115.57. SPEP with IgG kappa, M-spike  1.0; previous monoclonal kappa free light chain not detected.`"
  ],
  `"source_document`": `"2024-04-17_00:00:00.000_Progress_Notes_91591`"
}"
